$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the block of rows that hold the period "2508" data (rows 34-44, 11 workers)
# and insert it right after, creating rows 45-55 for the new period "2509".
$ws.Rows("34:44").Copy()
$ws.Rows("45:55").Insert()

# Update the "Periodo Mora" (column E) of the newly inserted rows from 2508 -> 2509
for ($r = 45; $r -le 55; $r++) {
    $ws.Cells.Item($r, 5).Value = "2509"
}

# Update the summary header: Cant. Periodos (F13) goes from 3 to 4
$ws.Range("F13").Value = 4

# Update the summary header: Valor Mora total (E11) goes from 1412112 to 2038452
$ws.Range("E11").Value = 2038452
